$wb = $excel.ActiveWorkbook

# Sheets that contain the duplicated event-listing table: "展览" and "全部类型"
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 231
    $ws.Range("F3").Value = 165
}
